$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing backlog items (text + estimation/priority numbers)
$ws.Range("A15").Value = "Game music"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 4

$ws.Range("A16").Value = "Additional players to choose from"
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = 2

# New backlog item added to the previously-blank row 17
$ws.Range("A17").Value = "Additional enemies types that can spawn"
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 2

# Extend the Total formula to include the new row
$ws.Range("B18").Formula = "=SUM(B3:B17)"

# Match the final selection left by the editor
$ws.Range("A14").Select()
